$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D18", "D19", "D20", "D22", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($cellref in $textCells) { $ws.Range($cellref).NumberFormat = "@" }

$ws.Range("D2").Value = "28.775.33"
$ws.Range("E2").Value = "  +2.69%  "

$ws.Range("D3").Value = "1.874.96"
$ws.Range("E3").Value = "  +2.41%  "

$ws.Range("D5").Value = "324.61"

$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("D7").Value = "0.4603"
$ws.Range("E7").Value = "  -0.97%  "

$ws.Range("D8").Value = "0.3873"
$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("D9").Value = "0.07863"
$ws.Range("E9").Value = "  +0.33%  "

$ws.Range("D10").Value = "0.9887"
$ws.Range("E10").Value = "  +3.13%  "

$ws.Range("D11").Value = "21.81"
$ws.Range("E11").Value = "  -0.35%  "

$ws.Range("D12").Value = "1.867.39"
$ws.Range("E12").Value = "  +2.32%  "

$ws.Range("D13").Value = "6.996"
$ws.Range("E13").Value = "  +1.56%  "

$ws.Range("D14").Value = "5.714"
$ws.Range("E14").Value = "  +0.61%  "

$ws.Range("D15").Value = "0.06986"
$ws.Range("E15").Value = "  +1.75%  "

$ws.Range("D16").Value = "88.49"
$ws.Range("E16").Value = "  +0.26%  "

$ws.Range("E17").Value = "  +0.26%  "

$ws.Range("D18").Value = "0.00001005"
$ws.Range("E18").Value = "  +1.35%  "

$ws.Range("D19").Value = "16.80"
$ws.Range("E19").Value = "  +1.01%  "

$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.26%  "

$ws.Range("D21").Value = "28.761.40"
$ws.Range("E21").Value = "  +2.61%  "

$ws.Range("D22").Value = "5.287"
$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("E23").Value = "  +0.84%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "2.100"
$ws.Range("E24").Value = "  +0.57%  "

$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "2.083.76"
$ws.Range("E25").Value = "  +1.97%  "

$ws.Range("D26").Value = "152.94"
$ws.Range("E26").Value = "  -1.23%  "

$ws.Range("D27").Value = "19.28"
$ws.Range("E27").Value = "  +0.81%  "

$ws.Range("D28").Value = "5.869"
$ws.Range("E28").Value = "  +3.88%  "

$ws.Range("D29").Value = "1.977"
$ws.Range("E29").Value = "  +0.95%  "

$ws.Range("D30").Value = "119.07"
$ws.Range("E30").Value = "  +0.63%  "

$ws.Range("D31").Value = "0.09322"
$ws.Range("E31").Value = "  +0.80%  "

$ws.Range("D32").Value = "0.9214"
$ws.Range("E32").Value = "  -1.34%  "

$ws.Range("D33").Value = "5.312"
$ws.Range("E33").Value = "  +1.11%  "

$ws.Range("D34").Value = "1.341"
$ws.Range("E34").Value = "  +1.68%  "

$ws.Range("D35").Value = "3.323"

$ws.Range("D36").Value = "0.05780"
$ws.Range("E36").Value = "  -1.09%  "

$ws.Range("D37").Value = "1.151"

$ws.Range("D38").Value = "0.02074"
$ws.Range("E38").Value = "  -2.17%  "

$ws.Range("D39").Value = "7.682"
$ws.Range("E39").Value = "  -0.66%  "

$ws.Range("D40").Value = "0.5645"
$ws.Range("E40").Value = "  +1.07%  "

$ws.Range("D41").Value = "0.1787"
$ws.Range("E41").Value = "  +1.73%  "

$ws.Range("D42").Value = "9.846"
$ws.Range("E42").Value = "  -0.27%  "

$ws.Range("D43").Value = "0.07214"
$ws.Range("E43").Value = "  -0.54%  "

$ws.Range("D44").Value = "11.73"
$ws.Range("E44").Value = "  +1.17%  "

$ws.Range("D45").Value = "0.5297"
$ws.Range("E45").Value = "  +0.73%  "

$ws.Range("D46").Value = "2.144"
$ws.Range("E46").Value = "  +2.45%  "

$ws.Range("D47").Value = "1.128"
$ws.Range("E47").Value = "  -1.53%  "

$ws.Range("E48").Value = "  +0.58%  "

$ws.Range("D49").Value = "113.54"
$ws.Range("E49").Value = "  +0.58%  "

$ws.Range("D50").Value = "2.416"
$ws.Range("E50").Value = "  +4.03%  "

$ws.Range("D51").Value = "1.004"
$ws.Range("E51").Value = "  +0.26%  "
